# Updates the Efna5-Epha1 LR-pairs sheet with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2457683333333333
$ws.Range("H2").Value = 0.737305
$ws.Range("I2").Value = 0.1447271191911903
$ws.Range("J2").Value = 0.1575855905380038
$ws.Range("M2").Value = 3.310473333333333
$ws.Range("N2").Value = 9.931419999999999
$ws.Range("O2").Value = 0.1683295705132556
$ws.Range("P2").Value = 0.1783110568845311
$ws.Range("Q2").Value = 0.8136095136777777
$ws.Range("R2").Value = 7.3224856231
$ws.Range("S2").Value = 0.02436185381507382
$ws.Range("T2").Value = 0.02809925319860441

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2457683333333333
$ws.Range("H3").Value = 0.737305
$ws.Range("I3").Value = 0.1447271191911903
$ws.Range("J3").Value = 0.1575855905380038
$ws.Range("M3").Value = 3.966196333333333
$ws.Range("O3").Value = 0.2016715007605908
$ws.Range("P3").Value = 0.2136300730433972
$ws.Range("Q3").Value = 0.9747654625161111
$ws.Range("R3").Value = 8.772889162644999
$ws.Range("S3").Value = 0.02918733532804426
$ws.Range("T3").Value = 0.03366502121722063

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2457683333333333
$ws.Range("H4").Value = 0.737305
$ws.Range("I4").Value = 0.1447271191911903
$ws.Range("J4").Value = 0.1575855905380038
$ws.Range("M4").Value = 5.013056
$ws.Range("N4").Value = 15.039168
$ws.Range("O4").Value = 0.2549017854764673
$ws.Range("P4").Value = 0.2700167690767302
$ws.Range("Q4").Value = 1.232050418026667
$ws.Range("R4").Value = 11.08845376224
$ws.Range("S4").Value = 0.03689120108869991
$ws.Range("T4").Value = 0.04255075201012031

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2457683333333333
$ws.Range("H5").Value = 0.737305
$ws.Range("I5").Value = 0.1447271191911903
$ws.Range("J5").Value = 0.1575855905380038
$ws.Range("M5").Value = 3.3026905
$ws.Range("N5").Value = 6.605381
$ws.Range("O5").Value = 0.1679338322424817
$ws.Range("P5").Value = 0.1185945682727144
$ws.Range("Q5").Value = 0.8116967397008333
$ws.Range("R5").Value = 4.870180438205
$ws.Range("S5").Value = 0.02430457975519102
$ws.Range("T5").Value = 0.0186887950758553

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2457683333333333
$ws.Range("H6").Value = 0.737305
$ws.Range("I6").Value = 0.1447271191911903
$ws.Range("J6").Value = 0.1575855905380038
$ws.Range("M6").Value = 4.074201666666667
$ws.Range("N6").Value = 12.222605
$ws.Range("O6").Value = 0.2071633110072045
$ws.Range("P6").Value = 0.2194475327226272
$ws.Range("Q6").Value = 1.001309753280556
$ws.Range("R6").Value = 9.011787779525001
$ws.Range("S6").Value = 0.02998214920418131
$ws.Range("T6").Value = 0.03458176903620312

# Row 7
$ws.Range("I7").Value = 0.6104826439049008
$ws.Range("J7").Value = 0.6647217777192627
$ws.Range("M7").Value = 3.310473333333333
$ws.Range("N7").Value = 9.931419999999999
$ws.Range("O7").Value = 0.1683295705132556
$ws.Range("P7").Value = 0.1783110568845311
$ws.Range("Q7").Value = 3.431937910406667
$ws.Range("R7").Value = 30.88744119366
$ws.Range("S7").Value = 0.1027622812543087
$ws.Range("T7").Value = 0.1185272427192861

# Row 8
$ws.Range("I8").Value = 0.6104826439049008
$ws.Range("J8").Value = 0.6647217777192627
$ws.Range("M8").Value = 3.966196333333333
$ws.Range("O8").Value = 0.2016715007605908
$ws.Range("P8").Value = 0.2136300730433972
$ws.Range("Q8").Value = 4.111720042999666
$ws.Range("R8").Value = 37.00548038699699
$ws.Range("S8").Value = 0.1231169509845947
$ws.Range("T8").Value = 0.1420045619277029

# Row 9
$ws.Range("I9").Value = 0.6104826439049008
$ws.Range("J9").Value = 0.6647217777192627
$ws.Range("M9").Value = 5.013056
$ws.Range("N9").Value = 15.039168
$ws.Range("O9").Value = 0.2549017854764673
$ws.Range("P9").Value = 0.2700167690767302
$ws.Range("Q9").Value = 5.196990037696
$ws.Range("R9").Value = 46.772910339264
$ws.Range("S9").Value = 0.1556131159337536
$ws.Range("T9").Value = 0.1794860267546957

# Row 10
$ws.Range("I10").Value = 0.6104826439049008
$ws.Range("J10").Value = 0.6647217777192627
$ws.Range("M10").Value = 3.3026905
$ws.Range("N10").Value = 6.605381
$ws.Range("O10").Value = 0.1679338322424817
$ws.Range("P10").Value = 0.1185945682727144
$ws.Range("Q10").Value = 3.4238695171355
$ws.Range("R10").Value = 20.543217102813
$ws.Range("S10").Value = 0.1025206899084723
$ws.Range("T10").Value = 0.07883239225008716

# Row 11
$ws.Range("I11").Value = 0.6104826439049008
$ws.Range("J11").Value = 0.6647217777192627
$ws.Range("M11").Value = 4.074201666666667
$ws.Range("N11").Value = 12.222605
$ws.Range("O11").Value = 0.2071633110072045
$ws.Range("P11").Value = 0.2194475327226272
$ws.Range("Q11").Value = 4.223688200018334
$ws.Range("R11").Value = 38.01319380016501
$ws.Range("S11").Value = 0.1264696058237714
$ws.Range("T11").Value = 0.1458715540674908

# Row 12
$ws.Range("G12").Value = 0.4156905
$ws.Range("H12").Value = 0.831381
$ws.Range("I12").Value = 0.2447902369039089
$ws.Range("J12").Value = 0.1776926317427335
$ws.Range("M12").Value = 3.310473333333333
$ws.Range("N12").Value = 9.931419999999999
$ws.Range("O12").Value = 0.1683295705132556
$ws.Range("P12").Value = 0.1783110568845311
$ws.Range("Q12").Value = 1.37613231517
$ws.Range("R12").Value = 8.256793891019999
$ws.Range("S12").Value = 0.04120543544387307
$ws.Range("T12").Value = 0.03168456096664058

# Row 13
$ws.Range("G13").Value = 0.4156905
$ws.Range("H13").Value = 0.831381
$ws.Range("I13").Value = 0.2447902369039089
$ws.Range("J13").Value = 0.1776926317427335
$ws.Range("M13").Value = 3.966196333333333
$ws.Range("O13").Value = 0.2016715007605908
$ws.Range("P13").Value = 0.2136300730433972
$ws.Range("Q13").Value = 1.6487101369015
$ws.Range("R13").Value = 9.892260821409
$ws.Range("S13").Value = 0.04936721444795187
$ws.Range("T13").Value = 0.03796048989847364

# Row 14
$ws.Range("G14").Value = 0.4156905
$ws.Range("H14").Value = 0.831381
$ws.Range("I14").Value = 0.2447902369039089
$ws.Range("J14").Value = 0.1776926317427335
$ws.Range("M14").Value = 5.013056
$ws.Range("N14").Value = 15.039168
$ws.Range("O14").Value = 0.2549017854764673
$ws.Range("P14").Value = 0.2700167690767302
$ws.Range("Q14").Value = 2.083879755168
$ws.Range("R14").Value = 12.503278531008
$ws.Range("S14").Value = 0.06239746845401378
$ws.Range("T14").Value = 0.04797999031191412

# Row 15
$ws.Range("G15").Value = 0.4156905
$ws.Range("H15").Value = 0.831381
$ws.Range("I15").Value = 0.2447902369039089
$ws.Range("J15").Value = 0.1776926317427335
$ws.Range("M15").Value = 3.3026905
$ws.Range("N15").Value = 6.605381
$ws.Range("O15").Value = 0.1679338322424817
$ws.Range("P15").Value = 0.1185945682727144
$ws.Range("Q15").Value = 1.37289706529025
$ws.Range("R15").Value = 5.491588261161
$ws.Range("S15").Value = 0.04110856257881839
$ws.Range("T15").Value = 0.0210733809467719

# Row 16
$ws.Range("G16").Value = 0.4156905
$ws.Range("H16").Value = 0.831381
$ws.Range("I16").Value = 0.2447902369039089
$ws.Range("J16").Value = 0.1776926317427335
$ws.Range("M16").Value = 4.074201666666667
$ws.Range("N16").Value = 12.222605
$ws.Range("O16").Value = 0.2071633110072045
$ws.Range("P16").Value = 0.2194475327226272
$ws.Range("Q16").Value = 1.6936069279175
$ws.Range("R16").Value = 10.161641567505
$ws.Range("S16").Value = 0.05071155597925173
$ws.Range("T16").Value = 0.03458176903620312
